# Apply numeric updates to Kujata_Profits leve-crafting sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2892.25
$ws.Range("I32").Value = 569
$ws.Range("J32").Value = 3666.6667
$ws.Range("K32").Value = 569
$ws.Range("L32").Value = 3666.6667
$ws.Range("M32").Value = -243
$ws.Range("N32").Value = -4318.6667

$ws.Range("H33").Value = 595.5714
$ws.Range("I33").Value = 618.1
$ws.Range("J33").Value = 539.25
$ws.Range("K33").Value = 618.1
$ws.Range("L33").Value = 539.25
$ws.Range("M33").Value = -389.1
$ws.Range("N33").Value = -997.25

$ws.Range("H40").Value = 1985.9546
$ws.Range("J40").Value = 2022.2307
$ws.Range("L40").Value = 2022.2307
$ws.Range("N40").Value = -2372.2307

$ws.Range("H70").Value = 1740.75
$ws.Range("I70").Value = 1705.7778
$ws.Range("J70").Value = 1785.7142
$ws.Range("K70").Value = 5117.3334
$ws.Range("L70").Value = 5357.142599999999
$ws.Range("M70").Value = -4847.3334
$ws.Range("N70").Value = -5897.142599999999

$ws.Range("H73").Value = 1740.75
$ws.Range("I73").Value = 1705.7778
$ws.Range("J73").Value = 1785.7142
$ws.Range("K73").Value = 5117.3334
$ws.Range("L73").Value = 5357.142599999999
$ws.Range("M73").Value = -4181.3334
$ws.Range("N73").Value = -7229.142599999999

$ws.Range("H111").Value = 6419.8887
$ws.Range("I111").Value = 7000
$ws.Range("J111").Value = 6129.8335
$ws.Range("K111").Value = 21000
$ws.Range("L111").Value = 18389.5005
$ws.Range("M111").Value = -17933
$ws.Range("N111").Value = -24523.5005

$ws.Range("H132").Value = 10108942
$ws.Range("I132").Value = 15159152
$ws.Range("J132").Value = 8520.817999999999
$ws.Range("K132").Value = 45477456
$ws.Range("L132").Value = 25562.454
$ws.Range("M132").Value = -45474926
$ws.Range("N132").Value = -30622.454

$ws.Range("H137").Value = 1537.4286
$ws.Range("I137").Value = 1075.0834
$ws.Range("J137").Value = 1884.1875
$ws.Range("K137").Value = 3225.2502
$ws.Range("L137").Value = 5652.5625
$ws.Range("M137").Value = -675.2501999999999
$ws.Range("N137").Value = -10752.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4252.6704
$ws.Range("I32").Value = 4152.171
$ws.Range("J32").Value = 6999.6665
$ws.Range("K32").Value = 4152.171
$ws.Range("L32").Value = 6999.6665
$ws.Range("M32").Value = -3865.171
$ws.Range("N32").Value = -7573.6665

$ws.Range("H45").Value = 1268
$ws.Range("I45").Value = 1087.5
$ws.Range("K45").Value = 1087.5
$ws.Range("M45").Value = -710.5

$ws.Range("H122").Value = 1733.7059
$ws.Range("I122").Value = 1301.1538
$ws.Range("K122").Value = 3903.4614
$ws.Range("M122").Value = -1453.4614

$ws.Range("H132").Value = 1750.0869
$ws.Range("I132").Value = 1282.5143
$ws.Range("J132").Value = 3237.818
$ws.Range("K132").Value = 3847.5429
$ws.Range("L132").Value = 9713.454000000002
$ws.Range("M132").Value = -1317.5429
$ws.Range("N132").Value = -14773.454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1486.4286
$ws.Range("J134").Value = 2249.25
$ws.Range("L134").Value = 6747.75
$ws.Range("N134").Value = -11817.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1358.683
$ws.Range("I31").Value = 1323.9
$ws.Range("K31").Value = 1323.9
$ws.Range("M31").Value = -1028.9

$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

$ws.Range("H34").Value = 1358.683
$ws.Range("I34").Value = 1323.9
$ws.Range("K34").Value = 1323.9
$ws.Range("M34").Value = -1121.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9317.5
$ws.Range("I3").Value = 6445.5557
$ws.Range("J3").Value = 13010
$ws.Range("K3").Value = 19336.6671
$ws.Range("L3").Value = 39030
$ws.Range("M3").Value = -19224.6671
$ws.Range("N3").Value = -39254

$ws.Range("H55").Value = 2000.6
$ws.Range("I55").Value = 454
$ws.Range("J55").Value = 2172.4443
$ws.Range("K55").Value = 1362
$ws.Range("L55").Value = 6517.3329
$ws.Range("M55").Value = -1185
$ws.Range("N55").Value = -6871.3329

$ws.Range("H92").Value = 229.37837
$ws.Range("I92").Value = 223.72414
$ws.Range("J92").Value = 249.875
$ws.Range("K92").Value = 671.17242
$ws.Range("L92").Value = 749.625
$ws.Range("M92").Value = 576.82758
$ws.Range("N92").Value = -3245.625

$ws.Range("H94").Value = 4109.467
$ws.Range("I94").Value = 3524
$ws.Range("J94").Value = 4199.5386
$ws.Range("K94").Value = 10572
$ws.Range("L94").Value = 12598.6158
$ws.Range("M94").Value = -9896
$ws.Range("N94").Value = -13950.6158

$ws.Range("H131").Value = 15873923
$ws.Range("J131").Value = 1045.8085
$ws.Range("L131").Value = 3137.4255
$ws.Range("N131").Value = -13217.4255

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 40000
$ws.Range("J51").Value = 40000
$ws.Range("L51").Value = 40000
$ws.Range("N51").Value = -41018

$ws.Range("H52").Value = 22999.5
$ws.Range("J52").Value = 22999.5
$ws.Range("L52").Value = 22999.5
$ws.Range("N52").Value = -23517.5

$ws.Range("H70").Value = 28128730
$ws.Range("I70").Value = 25003818
$ws.Range("J70").Value = 33336918
$ws.Range("K70").Value = 25003818
$ws.Range("L70").Value = 33336918
$ws.Range("M70").Value = -25003548
$ws.Range("N70").Value = -33337458

$ws.Range("H73").Value = 28128730
$ws.Range("I73").Value = 25003818
$ws.Range("J73").Value = 33336918
$ws.Range("K73").Value = 25003818
$ws.Range("L73").Value = 33336918
$ws.Range("M73").Value = -25002882
$ws.Range("N73").Value = -33338790

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H132").Value = 2973.4
$ws.Range("I132").Value = 2815.389
$ws.Range("J132").Value = 3210.4167
$ws.Range("K132").Value = 8446.167000000001
$ws.Range("L132").Value = 9631.250100000001
$ws.Range("M132").Value = -5916.167000000001
$ws.Range("N132").Value = -14691.2501

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2250
$ws.Range("J22").Value = 2250
$ws.Range("L22").Value = 2250
$ws.Range("N22").Value = -2840

$ws.Range("H27").Value = 2250
$ws.Range("J27").Value = 2250
$ws.Range("L27").Value = 2250
$ws.Range("N27").Value = -2464

$ws.Range("H46").Value = 5783.3335
$ws.Range("J46").Value = 6800
$ws.Range("L46").Value = 6800
$ws.Range("N46").Value = -7176

$ws.Range("H122").Value = 2579
$ws.Range("I122").Value = 2617.25
$ws.Range("K122").Value = 7851.75
$ws.Range("M122").Value = -5401.75

$ws.Range("H136").Value = 1697.1428
$ws.Range("I136").Value = 1176
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 3528
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -978
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2303.0476
$ws.Range("I132").Value = 1964.4
$ws.Range("J132").Value = 3149.6667
$ws.Range("K132").Value = 5893.200000000001
$ws.Range("L132").Value = 9449.000100000001
$ws.Range("M132").Value = -3363.200000000001
$ws.Range("N132").Value = -14509.0001

$ws.Range("H136").Value = 1141.5
$ws.Range("I136").Value = 998.29034
$ws.Range("J136").Value = 1483
$ws.Range("K136").Value = 2994.87102
$ws.Range("L136").Value = 4449
$ws.Range("M136").Value = -444.87102
$ws.Range("N136").Value = -9549
